$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-77 down to 16-78
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with its data.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are constant across all data rows in this sheet,
# so copy them from the row directly below (row 16, the row that used to be row 15).
# Note: use Value2 (not Value) for reading back cell contents in this runtime.
$ws.Range("A15").Value2 = $ws.Range("A16").Value2
$ws.Range("B15").Value2 = $ws.Range("B16").Value2
$ws.Range("C15").Value2 = $ws.Range("C16").Value2
$ws.Range("D15").Value2 = 44672
$ws.Range("E15").Value2 = $ws.Range("E16").Value2
$ws.Range("F15").Value2 = $ws.Range("F16").Value2
$ws.Range("G15").Value2 = $ws.Range("G16").Value2
$ws.Range("H15").Value2 = $ws.Range("H16").Value2
$ws.Range("I15").Value2 = $ws.Range("I16").Value2
$ws.Range("J15").Value2 = $ws.Range("J16").Value2
$ws.Range("K15").Value2 = $ws.Range("K16").Value2
$ws.Range("L15").Value2 = $ws.Range("L16").Value2
$ws.Range("M15").Value2 = 40
$ws.Range("N15").Value2 = 7000
$ws.Range("O15").Value2 = 8000
$ws.Range("P15").Value2 = 7500
$ws.Range("Q15").Value2 = $ws.Range("Q16").Value2
$ws.Range("R15").Value2 = "Perú"
$ws.Range("S15").Value2 = 1875
$ws.Range("T15").Value2 = $ws.Range("T16").Value2
